$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 "Save" - copy formatting (bold font, border, centered
# alignment) from the neighboring header cell G1, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# New data column values for rows 2 and 3 (both 0).
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
